$p = $ppt.ActivePresentation

function Get-ShapeById($container, $id) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# --- 1) Bump the "datetimeFigureOut" footer field cache (9/18/2020 -> 9/20/2020)
#        on the slide master and every slide layout.
function Update-DateField($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "9/18/2020") {
                $shp.TextFrame.TextRange.Text = "9/20/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateField $layouts.Item($li)
}

# --- 2) Slide 1: add a "last updated" textbox under the title.
$s1 = $p.Slides.Item(1)
$tb = $s1.Shapes.AddTextbox(1, 67.41803741455078, 135.33047485351562, 125.11000061035156, 29.081260681152344)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0
$tb.TextFrame.TextRange.Text = "20 Sept 2020"
$tb.Left = 67.41803741455078
$tb.Top = 135.33047485351562
$tb.Width = 125.11000061035156
$tb.Height = 29.081260681152344

# --- 3) Slide 2: swap the Red/Green config-button labels (and nudge/resize
#        the textbox callouts that now point at the other color).
$s2 = $p.Slides.Item(2)

$shp = Get-ShapeById $s2 24
$shp.Left = 267.23583984375
$shp.Top = 320.8048095703125
$shp.Width = 39.65858459472656
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ R/Y"

$shp = Get-ShapeById $s2 47
$shp.Left = 314.418212890625
$shp.Top = 431.0253601074219
$shp.Width = 52.65929412841797
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ Green"

$shp = Get-ShapeById $s2 52
$shp.Left = 455.8115234375
$shp.Top = 173.63204956054688
$shp.Width = 47.95834732055664
$shp.Height = 43.621891021728516
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ Green"

$shp = Get-ShapeById $s2 66
$shp.Left = 481.62158203125
$shp.Top = 272.7354431152344
$shp.Width = 46.167877197265625
$shp.Height = 43.621891021728516
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ Green"

$shp = Get-ShapeById $s2 67
$shp.Left = 466.669921875
$shp.Top = 480.8827819824219
$shp.Width = 52.04905700683594
$shp.Height = 43.621891021728516
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ Green"

$shp = Get-ShapeById $s2 72
$shp.Left = 623.05908203125
$shp.Top = 171.36378479003906
$shp.Width = 47.95834732055664
$shp.Height = 43.621891021728516
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ Green"

$shp = Get-ShapeById $s2 77
$shp.Left = 404.5920715332031
$shp.Top = 325.86260986328125
$shp.Width = 39.65858459472656
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ R/Y"

$shp = Get-ShapeById $s2 78
$shp.Left = 432.7739562988281
$shp.Top = 379.2055358886719
$shp.Width = 39.65858459472656
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ R/Y"

$shp = Get-ShapeById $s2 79
$shp.Left = 419.95220947265625
$shp.Top = 431.39874267578125
$shp.Width = 39.65858459472656
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ R/Y"

$shp = Get-ShapeById $s2 96
$shp.Left = 609.8745727539062
$shp.Top = 285.2762451171875
$shp.Width = 44.938899993896484
$shp.Height = 43.621891021728516
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ Green"

$shp = Get-ShapeById $s2 101
$shp.Left = 730.0072631835938
$shp.Top = 217.58425903320312
$shp.Width = 39.65858459472656
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ R/Y"

$shp = Get-ShapeById $s2 102
$shp.Left = 741.414306640625
$shp.Top = 325.86260986328125
$shp.Width = 39.65858459472656
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ R/Y"

$shp = Get-ShapeById $s2 103
$shp.Left = 742.5516967773438
$shp.Top = 434.1409606933594
$shp.Width = 39.65858459472656
$shp.Height = 31.504724502563477
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ R/Y"

$shp = Get-ShapeById $s2 104
$shp.Left = 631.7423095703125
$shp.Top = 486.2687683105469
$shp.Width = 49.73409652709961
$shp.Height = 43.621891021728516
$shp.TextFrame.TextRange.Paragraphs(2).Text = "zzz"
$shp.TextFrame.TextRange.Paragraphs(2).Text = "+ Green"

